$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price / volume snapshot (GitHub Actions scheduled update).
# D (Price) values are written with a leading apostrophe so Excel keeps them
# as literal text (matching the workbook's inlineStr storage) instead of
# auto-converting number-looking strings like "1.022" or "13.00" into reals.

$ws.Range("D2").Value2 = "'28.408.88"
$ws.Range("E2").Value2 = "  +0.41%  "
$ws.Range("D3").Value2 = "'1.872.63"
$ws.Range("E3").Value2 = "  -0.80%  "
$ws.Range("D4").Value2 = "'1.022"
$ws.Range("E4").Value2 = "  +1.16%  "
$ws.Range("D5").Value2 = "'316.89"
$ws.Range("E5").Value2 = "  +0.16%  "
$ws.Range("D6").Value2 = "'1.019"
$ws.Range("E6").Value2 = "  +0.86%  "
$ws.Range("D7").Value2 = "'0.5108"
$ws.Range("E7").Value2 = "  -0.74%  "
$ws.Range("E8").Value2 = "  +1.09%  "
$ws.Range("D9").Value2 = "'0.08452"
$ws.Range("E9").Value2 = "  +0.64%  "
$ws.Range("D10").Value2 = "'1.109"
$ws.Range("E10").Value2 = "  -1.45%  "
$ws.Range("D11").Value2 = "'41.98"
$ws.Range("E11").Value2 = "  +0.58%  "
$ws.Range("D12").Value2 = "'6.251"
$ws.Range("D13").Value2 = "'1.874.91"
$ws.Range("E13").Value2 = "  -0.29%  "
$ws.Range("D14").Value2 = "'20.46"
$ws.Range("E14").Value2 = "  -0.78%  "
$ws.Range("D17").Value2 = "'0.00001111"
$ws.Range("E17").Value2 = "  +0.31%  "
$ws.Range("D18").Value2 = "'90.99"
$ws.Range("E18").Value2 = "  -0.17%  "
$ws.Range("D19").Value2 = "'0.06773"
$ws.Range("E19").Value2 = "  +1.13%  "
$ws.Range("D20").Value2 = "'17.73"
$ws.Range("E20").Value2 = "  -0.57%  "
$ws.Range("D21").Value2 = "'1.019"
$ws.Range("E21").Value2 = "  +0.96%  "
$ws.Range("E22").Value2 = "  -1.74%  "
$ws.Range("D23").Value2 = "'28.467.31"
$ws.Range("E23").Value2 = "  +0.53%  "
$ws.Range("D24").Value2 = "'11.17"
$ws.Range("E24").Value2 = "  -0.10%  "
$ws.Range("D25").Value2 = "'2.293"
$ws.Range("E25").Value2 = "  -0.02%  "
$ws.Range("D26").Value2 = "'2.079.77"
$ws.Range("E26").Value2 = "  -0.86%  "
$ws.Range("D27").Value2 = "'161.71"
$ws.Range("E27").Value2 = "  +0.68%  "
$ws.Range("D28").Value2 = "'20.76"
$ws.Range("E28").Value2 = "  +0.14%  "
$ws.Range("D29").Value2 = "'2.351"
$ws.Range("E29").Value2 = "  -4.84%  "
$ws.Range("D30").Value2 = "'127.09"
$ws.Range("E30").Value2 = "  +1.11%  "
$ws.Range("D31").Value2 = "'0.1052"
$ws.Range("E31").Value2 = "  -0.82%  "
$ws.Range("E32").Value2 = "  -0.25%  "
$ws.Range("D33").Value2 = "'5.766"
$ws.Range("E33").Value2 = "  -2.22%  "
$ws.Range("D34").Value2 = "'3.635"
$ws.Range("E34").Value2 = "  +0.23%  "
$ws.Range("D35").Value2 = "'0.02433"
$ws.Range("E35").Value2 = "  -0.66%  "
$ws.Range("D36").Value2 = "'0.06458"
$ws.Range("E36").Value2 = "  -2.02%  "
$ws.Range("D37").Value2 = "'0.2175"
$ws.Range("E37").Value2 = "  -1.73%  "
$ws.Range("D38").Value2 = "'8.835"
$ws.Range("E38").Value2 = "  -6.90%  "
$ws.Range("E39").Value2 = "  +0.94%  "
$ws.Range("E40").Value2 = "  -1.56%  "
$ws.Range("D41").Value2 = "'0.6365"
$ws.Range("E41").Value2 = "  -2.23%  "
$ws.Range("D42").Value2 = "'4.980"
$ws.Range("E42").Value2 = "  -0.75%  "
$ws.Range("E43").Value2 = "  -0.40%  "
$ws.Range("D44").Value2 = "'0.6022"
$ws.Range("E44").Value2 = "  -1.46%  "
$ws.Range("D45").Value2 = "'13.00"
$ws.Range("E45").Value2 = "  -0.67%  "
$ws.Range("D46").Value2 = "'3.710"
$ws.Range("E46").Value2 = "  +0.24%  "
$ws.Range("D49").Value2 = "'1.206"
$ws.Range("E49").Value2 = "  -2.63%  "
$ws.Range("D50").Value2 = "'121.74"
$ws.Range("E50").Value2 = "  +0.34%  "
$ws.Range("D51").Value2 = "'0.06844"
$ws.Range("E51").Value2 = "  -1.09%  "

# Rows 15/16 swapped order (Chainlink <-> BinanceUSD) with refreshed data
$ws.Range("B15").Value2 = "BinanceUSD"
$ws.Range("C15").Value2 = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D15").Value2 = "'1.022"
$ws.Range("E15").Value2 = "  +1.18%  "
$ws.Range("B16").Value2 = "Chainlink"
$ws.Range("C16").Value2 = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value2 = "'7.229"
$ws.Range("E16").Value2 = "  -0.69%  "

# Rows 47/48 swapped order (WEMIXTOKEN <-> NEARProtocol) with refreshed data
$ws.Range("B47").Value2 = "NEARProtocol"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value2 = "'1.991"
$ws.Range("E47").Value2 = "  -1.49%  "
$ws.Range("B48").Value2 = "WEMIXTOKEN"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value2 = "'1.208"
$ws.Range("E48").Value2 = "  -6.00%  "
